$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark from the end of the "Worked cross-
#        functionally ... recommendations" bullet onto the
#        "https://sshofa.github.io" link in the contact-info table. ---
#
# Real Word keeps exactly one "_GoBack" bookmark, silently relocating it to
# mark the position of the most recent edit. Remove it from its old spot
# and recreate it around the new text.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content.Duplicate
$target.Find.Execute("https://sshofa.github.io", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
